$wb = $excel.ActiveWorkbook

# --- Sheet "Obs": swap Crefbiascv / Brefbiascv text (A25 / A26) ---
$wsObs = $wb.Worksheets.Item("Obs")
$wsObs.Range("A25").Value = "Brefbiascv"
$wsObs.Range("A26").Value = "Crefbiascv"

# --- Sheet "Imp": swap TACSD/TACFrac, TAESD/TAEFrac, SizeLimSD/SizeLimFrac ---
$wsImp = $wb.Worksheets.Item("Imp")
$wsImp.Range("A3").Value = "TACFrac"
$wsImp.Range("A4").Value = "TACSD"
$wsImp.Range("A5").Value = "TAEFrac"
$wsImp.Range("A6").Value = "TAESD"
$wsImp.Range("A7").Value = "SizeLimFrac"
$wsImp.Range("A8").Value = "SizeLimSD"

# --- Sheet "Imp": remove row 9 (Source) ---
$wsImp.Rows.Item(9).Delete()
